$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '63.548.70'
Set-TextValue 'E2' '  -4.54%  '
Set-TextValue 'D3' '3.594.54'
Set-TextValue 'E3' '  +0.02%  '
Set-TextValue 'E4' '  +0.40%  '
Set-TextValue 'D5' '402.84'
Set-TextValue 'E5' '  -2.93%  '
Set-TextValue 'D6' '132.50'
Set-TextValue 'E6' '  +1.54%  '
Set-TextValue 'D7' '3.584.10'
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.611'
Set-TextValue 'E8' '  -5.92%  '
Set-TextValue 'E9' '  +0.21%  '
Set-TextValue 'D10' '0.714'
Set-TextValue 'E10' '  -7.48%  '
Set-TextValue 'D11' '0.154'
Set-TextValue 'E11' '  -11.12%  '
Set-TextValue 'D12' '0.0000304'
Set-TextValue 'E12' '  -6.92%  '
Set-TextValue 'D13' '41.22'
Set-TextValue 'E13' '  -2.80%  '
Set-TextValue 'D14' '9.75'
Set-TextValue 'E14' '  -1.21%  '
Set-TextValue 'D15' '4.179.04'
Set-TextValue 'E15' '  +0.36%  '
Set-TextValue 'E16' '  -1.27%  '
Set-TextValue 'D17' '3.598.21'
Set-TextValue 'E17' '  -0.65%  '
Set-TextValue 'D18' '19.66'
Set-TextValue 'E18' '  -3.26%  '
Set-TextValue 'D19' '13.22'
Set-TextValue 'E19' '  +7.27%  '
Set-TextValue 'E20' '  -6.49%  '
Set-TextValue 'D21' '63.810.64'
Set-TextValue 'E21' '  -4.09%  '
Set-TextValue 'D22' '414.31'
Set-TextValue 'E22' '  -7.37%  '
Set-TextValue 'D23' '14.85'
Set-TextValue 'E23' '  +13.64%  '
Set-TextValue 'D24' '84.29'
Set-TextValue 'E24' '  -5.19%  '
Set-TextValue 'D25' '2.95'
Set-TextValue 'E25' '  -6.33%  '
Set-TextValue 'D26' '35.05'
Set-TextValue 'E26' '  -1.08%  '
Set-TextValue 'D27' '3.14'
Set-TextValue 'E27' '  -5.46%  '
Set-TextValue 'D28' '9.28'
Set-TextValue 'E28' '  -6.74%  '
Set-TextValue 'E29' '  +5.55%  '
Set-TextValue 'D30' '12.57'
Set-TextValue 'E30' '  +1.96%  '
Set-TextValue 'D31' '2.69'
Set-TextValue 'E31' '  -2.75%  '
Set-TextValue 'D32' '0.114'
Set-TextValue 'E32' '  -2.38%  '
Set-TextValue 'D33' '6.85'
Set-TextValue 'E33' '  -7.33%  '
Set-TextValue 'E34' '  -1.37%  '
Set-TextValue 'D35' '40.40'
Set-TextValue 'E35' '  +1.40%  '
Set-TextValue 'E36' '  -0.11%  '
Set-TextValue 'D37' '55.17'
Set-TextValue 'E37' '  -2.60%  '
Set-TextValue 'D38' '0.0458'
Set-TextValue 'E38' '  -6.78%  '
Set-TextValue 'D39' '2.84'
Set-TextValue 'E39' '  +24.10%  '
Set-TextValue 'D40' '0.997'
Set-TextValue 'E40' '  -0.14%  '
Set-TextValue 'D41' '0.138'
Set-TextValue 'E41' '  -5.52%  '
Set-TextValue 'D42' '3.14'
Set-TextValue 'E42' '  +22.69%  '
Set-TextValue 'B43' 'Monero'
Set-TextValue 'C43' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D43' '144.10'
Set-TextValue 'E43' '  -3.42%  '
Set-TextValue 'B44' 'NEARProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D44' '4.35'
Set-TextValue 'E44' '  +0.94%  '
Set-TextValue 'D45' '0.0₃0626'
Set-TextValue 'E45' '  -12.44%  '
Set-TextValue 'B46' 'LidoDAOToken'
Set-TextValue 'C46' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D46' '3.25'
Set-TextValue 'E46' '  -0.18%  '
Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '26.11'
Set-TextValue 'E47' '  +22.16%  '
Set-TextValue 'D48' '2.04'
Set-TextValue 'E48' '  +3.23%  '
Set-TextValue 'D49' '2.79'
Set-TextValue 'E49' '  -6.90%  '
Set-TextValue 'D50' '2.52'
Set-TextValue 'E50' '  -7.71%  '
Set-TextValue 'D51' '0.286'
Set-TextValue 'E51' '  -8.00%  '
